$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 381.9091
$ws.Range("I12").Value = 322.33334
$ws.Range("K12").Value = 322.33334
$ws.Range("M12").Value = -152.33334
# Row 33
$ws.Range("H33").Value = 502.8889
$ws.Range("J33").Value = 1097
$ws.Range("L33").Value = 1097
$ws.Range("N33").Value = -1555
# Row 64
$ws.Range("H64").Value = 7661.8945
$ws.Range("I64").Value = 6397.3335
$ws.Range("J64").Value = 8800
$ws.Range("K64").Value = 6397.3335
$ws.Range("L64").Value = 8800
$ws.Range("M64").Value = -6149.3335
$ws.Range("N64").Value = -9296
# Row 67
$ws.Range("H67").Value = 7661.8945
$ws.Range("I67").Value = 6397.3335
$ws.Range("J67").Value = 8800
$ws.Range("K67").Value = 6397.3335
$ws.Range("L67").Value = 8800
$ws.Range("M67").Value = -5539.3335
$ws.Range("N67").Value = -10516
# Row 87
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
# Row 90
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
# Row 98
$ws.Range("H98").Value = 443.26666
$ws.Range("I98").Value = 443.26666
$ws.Range("K98").Value = 443.26666
$ws.Range("M98").Value = 1054.73334
# Row 122
$ws.Range("H122").Value = 443.26666
$ws.Range("I122").Value = 443.26666
$ws.Range("K122").Value = 1329.79998
$ws.Range("M122").Value = 1120.20002
# Row 127
$ws.Range("H127").Value = 9531.973
$ws.Range("I127").Value = 2182
$ws.Range("K127").Value = 6546
$ws.Range("M127").Value = -1586
# Row 137
$ws.Range("H137").Value = 2315.9412
$ws.Range("I137").Value = 1810.3334
$ws.Range("K137").Value = 5431.0002
$ws.Range("M137").Value = -2881.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13870.356
$ws.Range("I32").Value = 10119.698
$ws.Range("K32").Value = 10119.698
$ws.Range("M32").Value = -9832.698
# Row 45
$ws.Range("H45").Value = 5075.0625
$ws.Range("I45").Value = 2862.6667
$ws.Range("K45").Value = 2862.6667
$ws.Range("M45").Value = -2485.6667
# Row 132
$ws.Range("H132").Value = 4927.189
$ws.Range("I132").Value = 5082.4614
$ws.Range("J132").Value = 4560.1816
$ws.Range("K132").Value = 15247.3842
$ws.Range("L132").Value = 13680.5448
$ws.Range("M132").Value = -12717.3842
$ws.Range("N132").Value = -18740.5448

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2309.625
$ws.Range("I20").Value = 1948.75
$ws.Range("J20").Value = 2670.5
$ws.Range("K20").Value = 1948.75
$ws.Range("L20").Value = 2670.5
$ws.Range("M20").Value = -1701.75
$ws.Range("N20").Value = -3164.5
# Row 94
$ws.Range("H94").Value = 25004360
$ws.Range("J94").Value = 40005000
$ws.Range("L94").Value = 40005000
$ws.Range("N94").Value = -40005902
# Row 99
$ws.Range("H99").Value = 757.5454999999999
$ws.Range("I99").Value = 757.5454999999999
$ws.Range("K99").Value = 757.5454999999999
$ws.Range("M99").Value = 740.4545000000001
# Row 134
$ws.Range("H134").Value = 4448.857
$ws.Range("I134").Value = 2817.4614
$ws.Range("J134").Value = 7099.875
$ws.Range("K134").Value = 8452.3842
$ws.Range("L134").Value = 21299.625
$ws.Range("M134").Value = -5917.3842
$ws.Range("N134").Value = -26369.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 291.83334
$ws.Range("I5").Value = 424
$ws.Range("J5").Value = 265.4
$ws.Range("K5").Value = 424
$ws.Range("L5").Value = 265.4
$ws.Range("M5").Value = -312
$ws.Range("N5").Value = -489.4
# Row 8
$ws.Range("H8").Value = 800
$ws.Range("I8").Value = 800
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 800
$ws.Range("M8").Value = -660
$ws.Range("N8").ClearContents()
# Row 10
$ws.Range("H10").Value = 390.2
$ws.Range("I10").Value = 390.2
$ws.Range("K10").Value = 390.2
$ws.Range("M10").Value = -251.2
# Row 12
$ws.Range("H12").Value = 594.75
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 659.6667
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 659.6667
$ws.Range("M12").Value = -230
$ws.Range("N12").Value = -999.6667
# Row 15
$ws.Range("H15").Value = 745
$ws.Range("J15").Value = 745
$ws.Range("L15").Value = 745
$ws.Range("N15").Value = -1085
# Row 31
$ws.Range("H31").Value = 3127.1973
$ws.Range("I31").Value = 2136.825
$ws.Range("J31").Value = 4405.0967
$ws.Range("K31").Value = 2136.825
$ws.Range("L31").Value = 4405.0967
$ws.Range("M31").Value = -1841.825
$ws.Range("N31").Value = -4995.0967
# Row 34
$ws.Range("H34").Value = 3127.1973
$ws.Range("I34").Value = 2136.825
$ws.Range("J34").Value = 4405.0967
$ws.Range("K34").Value = 2136.825
$ws.Range("L34").Value = 4405.0967
$ws.Range("M34").Value = -1934.825
$ws.Range("N34").Value = -4809.0967

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 243.44444
$ws.Range("J12").Value = 240.625
$ws.Range("L12").Value = 721.875
$ws.Range("N12").Value = -1067.875
# Row 29
$ws.Range("H29").Value = 31624744
$ws.Range("I29").Value = 9723483
$ws.Range("K29").Value = 29170449
$ws.Range("M29").Value = -29170172
# Row 114
$ws.Range("H114").Value = 1494.7693
$ws.Range("I114").Value = 624
$ws.Range("K114").Value = 1872
$ws.Range("M114").Value = 1382
# Row 117
$ws.Range("H117").Value = 2187.1052
$ws.Range("I117").Value = 1009.8571
$ws.Range("J117").Value = 2873.8333
$ws.Range("K117").Value = 3029.5713
$ws.Range("L117").Value = 8621.499899999999
$ws.Range("M117").Value = 412.4287000000004
$ws.Range("N117").Value = -15505.4999
# Row 136
$ws.Range("H136").Value = 2995.158
$ws.Range("I136").Value = 1901.5454
$ws.Range("K136").Value = 5704.6362
$ws.Range("M136").Value = -604.6361999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3938.2942
$ws.Range("I126").Value = 1097.2
$ws.Range("K126").Value = 3291.6
$ws.Range("M126").Value = -821.6000000000004
# Row 132
$ws.Range("H132").Value = 4885.087
$ws.Range("I132").Value = 4054.0833
$ws.Range("J132").Value = 5791.636
$ws.Range("K132").Value = 12162.2499
$ws.Range("L132").Value = 17374.908
$ws.Range("M132").Value = -9632.249899999999
$ws.Range("N132").Value = -22434.908
# Row 140
$ws.Range("H140").Value = 79658.336
$ws.Range("J140").Value = 79658.336
$ws.Range("L140").Value = 79658.336
$ws.Range("N140").Value = -90018.336
# Row 141
$ws.Range("H141").Value = 124999.5
$ws.Range("J141").Value = 124999.5
$ws.Range("L141").Value = 124999.5
$ws.Range("N141").Value = -135359.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 162.9
$ws.Range("I16").Value = 162.9
$ws.Range("K16").Value = 162.9
$ws.Range("M16").Value = 7.099999999999994
# Row 22
$ws.Range("H22").Value = 2050
$ws.Range("I22").Value = 1775
$ws.Range("J22").Value = 2875
$ws.Range("K22").Value = 1775
$ws.Range("L22").Value = 2875
$ws.Range("M22").Value = -1480
$ws.Range("N22").Value = -3465
# Row 27
$ws.Range("H27").Value = 2050
$ws.Range("I27").Value = 1775
$ws.Range("J27").Value = 2875
$ws.Range("K27").Value = 1775
$ws.Range("L27").Value = 2875
$ws.Range("M27").Value = -1668
$ws.Range("N27").Value = -3089
# Row 40
$ws.Range("H40").Value = 13979.923
$ws.Range("I40").Value = 20905
$ws.Range("K40").Value = 20905
$ws.Range("M40").Value = -20769
# Row 46
$ws.Range("H46").Value = 9245.129000000001
$ws.Range("J46").Value = 10269.23
$ws.Range("L46").Value = 10269.23
$ws.Range("N46").Value = -10645.23
# Row 122
$ws.Range("H122").Value = 5405.5
$ws.Range("I122").Value = 4253
$ws.Range("K122").Value = 12759
$ws.Range("M122").Value = -10309
# Row 136
$ws.Range("H136").Value = 6115.067
$ws.Range("I136").Value = 3855.4
$ws.Range("J136").Value = 7244.9
$ws.Range("K136").Value = 11566.2
$ws.Range("L136").Value = 21734.7
$ws.Range("M136").Value = -9016.200000000001
$ws.Range("N136").Value = -26834.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 47623864
$ws.Range("I81").Value = 4633.6665
$ws.Range("K81").Value = 9267.333000000001
$ws.Range("M81").Value = -8206.333000000001
# Row 84
$ws.Range("H84").Value = 47623864
$ws.Range("I84").Value = 4633.6665
$ws.Range("K84").Value = 46336.665
$ws.Range("M84").Value = -41032.665
# Row 112
$ws.Range("H112").Value = 39903.227
$ws.Range("J112").Value = 39903.227
$ws.Range("L112").Value = 39903.227
$ws.Range("N112").Value = -42857.227
# Row 113
$ws.Range("H113").Value = 795.25
$ws.Range("I113").Value = 794.44446
$ws.Range("J113").Value = 797.6667
$ws.Range("K113").Value = 2383.33338
$ws.Range("L113").Value = 2393.0001
$ws.Range("M113").Value = -213.33338
$ws.Range("N113").Value = -6733.0001
# Row 122
$ws.Range("H122").Value = 3904.7368
$ws.Range("I122").Value = 2591.1538
$ws.Range("K122").Value = 7773.4614
$ws.Range("M122").Value = -5323.4614
# Row 126
$ws.Range("H126").Value = 1567.8462
$ws.Range("I126").Value = 1448.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4345.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1875.5
$ws.Range("N126").Value = -13940
# Row 132
$ws.Range("H132").Value = 3404.138
$ws.Range("I132").Value = 2119.3333
$ws.Range("K132").Value = 6357.999899999999
$ws.Range("M132").Value = -3827.999899999999
